$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a row for "id_membre" after the header (new row 3) ---
$ws.Rows.Item(3).Insert()

# --- 2. Insert three rows before the old "tel_membre" row (which is now
#        row 7 after the previous insert) to hold adresse_membre,
#        cp_membre and ville_membre ---
$ws.Rows.Item(7).Resize(3).Insert()

$ws.Application.CutCopyMode = $false

# ===================== ROW 3 : id_membre =====================
# Build the base formatting for the "plain" (border all around, no
# special alignment) cells by copying an existing data cell format.
$ws.Range("A4").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)

# Build the base formatting for the "bottom-less border" cells (used by
# columns F/G/H) by copying an existing header-row-2 cell format.
$ws.Range("F2").Copy()
$ws.Range("F3:H3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A3").Value = "id_membre"
$ws.Range("B3").Value = "Numérique"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = "×"
$ws.Range("G3").Value = "Auto incrémenté"

$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("E3").HorizontalAlignment = -4131
$ws.Range("E3").VerticalAlignment = -4108

$ws.Range("C3").HorizontalAlignment = -4152
$ws.Range("C3").VerticalAlignment = -4108

$ws.Range("F3").Font.Bold = $false
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").VerticalAlignment = -4108

$ws.Range("G3").Font.Bold = $false
$ws.Range("G3").HorizontalAlignment = -4131
$ws.Range("G3").VerticalAlignment = -4108

$ws.Range("H3").Font.Bold = $false
$ws.Range("H3").HorizontalAlignment = -4131
$ws.Range("H3").VerticalAlignment = -4108
$ws.Range("H3").WrapText = $true

# ===================== ROWS 4-6 (unchanged data, shifted) =====================
# email_membre / nom_membre / prenom_membre already kept their values and
# formatting after the row insert - nothing further to do.

# ===================== ROWS 7-9 : adresse_membre, cp_membre, ville_membre =====================
$ws.Range("A6").Copy()
$ws.Range("A7:H9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A7").Value = "adresse_membre"
$ws.Range("B7").Value = "Alphabétique"
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = "×"

$ws.Range("A8").Value = "cp_membre"
$ws.Range("B8").Value = "Alphanumérique"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "×"

$ws.Range("A9").Value = "ville_membre"
$ws.Range("B9").Value = "Alphabétique"
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = "×"

# ===================== ROWS 10-11 : tel_membre, mot_passe_membre =====================
$ws.Range("B10").Value = "Numérique"

# --- Selection / view state ---
$ws.Range("C10").Select()

Write-Host "done"
